$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 641, shifting existing rows 641:687 down to 642:688
$ws.Rows.Item(641).Insert()

# Populate the newly inserted row 641 with the new market-report entry
$ws.Cells.Item(641, 1).Value = 10
$ws.Cells.Item(641, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(641, 3).Value = "La Araucanía"
$ws.Cells.Item(641, 4).Value = 44461
$ws.Cells.Item(641, 5).Value = 9
$ws.Cells.Item(641, 6).Value = "Fruta"
$ws.Cells.Item(641, 7).Value = 100106
$ws.Cells.Item(641, 8).Value = "Oleaginosos"
$ws.Cells.Item(641, 9).Value = 100106002
$ws.Cells.Item(641, 10).Value = "Palta"
$ws.Cells.Item(641, 11).Value = "Hass"
$ws.Cells.Item(641, 12).Value = "1a nueva(o)"
$ws.Cells.Item(641, 13).Value = 80
$ws.Cells.Item(641, 14).Value = 3500
$ws.Cells.Item(641, 15).Value = 3500
$ws.Cells.Item(641, 16).Value = 3500
$ws.Cells.Item(641, 17).Value = "$/kilo (en bandeja de 18 kilos)"
$ws.Cells.Item(641, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(641, 19).Value = 3500
$ws.Cells.Item(641, 20).Value = 1
